$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-93 down to 22-94
$ws.Rows.Item(21).Insert()

# Fill in the new row 21 with data
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = 44592
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = 100112001
$ws.Range("G21").Value = "Berenjena"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 7000
$ws.Range("N21").Value = "`$/caja 50 unidades"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 140
$ws.Range("Q21").Value = 50
$ws.Range("R21").Value = "Hortaliza"
